$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 1952.6
$ws.Cells.Item(34, 9).Value = 1952.6
$ws.Cells.Item(34, 11).Value = 1952.6
$ws.Cells.Item(34, 13).Value = -1749.6
$ws.Cells.Item(36, 8).Value = 1952.6
$ws.Cells.Item(36, 9).Value = 1952.6
$ws.Cells.Item(36, 11).Value = 1952.6
$ws.Cells.Item(36, 13).Value = -1237.6
$ws.Cells.Item(40, 8).Value = 1761.2
$ws.Cells.Item(40, 9).Value = 1701
$ws.Cells.Item(40, 11).Value = 1701
$ws.Cells.Item(40, 13).Value = -1526
$ws.Cells.Item(64, 8).Value = 4000
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 4000
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 4000
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -4496
$ws.Cells.Item(67, 8).Value = 4000
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 4000
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 4000
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -5716
$ws.Cells.Item(82, 8).Value = 1800
$ws.Cells.Item(82, 9).Value = 1800
$ws.Cells.Item(82, 11).Value = 5400
$ws.Cells.Item(82, 13).Value = -4994
$ws.Cells.Item(85, 8).Value = 1800
$ws.Cells.Item(85, 9).Value = 1800
$ws.Cells.Item(85, 11).Value = 5400
$ws.Cells.Item(85, 13).Value = -3996
$ws.Cells.Item(103, 8).Value = 964
$ws.Cells.Item(103, 9).Value = 433.33334
$ws.Cells.Item(103, 11).Value = 1300.00002
$ws.Cells.Item(103, 13).Value = -714.00002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 4298.3335
$ws.Cells.Item(63, 9).Value = 2900
$ws.Cells.Item(63, 11).Value = 2900
$ws.Cells.Item(63, 13).Value = -2214
$ws.Cells.Item(66, 8).Value = 4298.3335
$ws.Cells.Item(66, 9).Value = 2900
$ws.Cells.Item(66, 11).Value = 14500
$ws.Cells.Item(66, 13).Value = -11068
$ws.Cells.Item(74, 8).Value = 2000
$ws.Cells.Item(74, 9).Value = 2000
$ws.Cells.Item(74, 11).Value = 2000
$ws.Cells.Item(74, 13).Value = -1126
$ws.Cells.Item(77, 8).Value = 2000
$ws.Cells.Item(77, 9).Value = 2000
$ws.Cells.Item(77, 11).Value = 10000
$ws.Cells.Item(77, 13).Value = -5632
$ws.Cells.Item(88, 8).Value = 1972.75
$ws.Cells.Item(88, 10).Value = 1997
$ws.Cells.Item(88, 12).Value = 1997
$ws.Cells.Item(88, 14).Value = -2809
$ws.Cells.Item(91, 8).Value = 1972.75
$ws.Cells.Item(91, 10).Value = 1997
$ws.Cells.Item(91, 12).Value = 1997
$ws.Cells.Item(91, 14).Value = -4805

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 13).ClearContents()
$ws.Cells.Item(29, 8).Value = 19500
$ws.Cells.Item(29, 9).Value = 19500
$ws.Cells.Item(29, 11).Value = 19500
$ws.Cells.Item(29, 13).Value = -19211
$ws.Cells.Item(36, 8).Value = 7000
$ws.Cells.Item(36, 10).Value = 1000
$ws.Cells.Item(36, 12).Value = 1000
$ws.Cells.Item(36, 14).Value = -2068
$ws.Cells.Item(94, 8).Value = 2801.8
$ws.Cells.Item(94, 9).Value = 1333
$ws.Cells.Item(94, 11).Value = 1333
$ws.Cells.Item(94, 13).Value = -882
$ws.Cells.Item(99, 8).Value = 1140.7142
$ws.Cells.Item(99, 9).Value = 1328.3334
$ws.Cells.Item(99, 11).Value = 1328.3334
$ws.Cells.Item(99, 13).Value = 169.6666
$ws.Cells.Item(105, 8).Value = 3699.8333
$ws.Cells.Item(105, 9).Value = 3199.75
$ws.Cells.Item(105, 10).Value = 4700
$ws.Cells.Item(105, 11).Value = 3199.75
$ws.Cells.Item(105, 12).Value = 4700
$ws.Cells.Item(105, 13).Value = -1452.75
$ws.Cells.Item(105, 14).Value = -8194

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 7591082
$ws.Cells.Item(11, 9).Value = 8500271
$ws.Cells.Item(11, 11).Value = 8500271
$ws.Cells.Item(11, 13).Value = -8500132
$ws.Cells.Item(70, 8).Value = 8387.833
$ws.Cells.Item(70, 9).Value = 7109
$ws.Cells.Item(70, 10).Value = 9666.667
$ws.Cells.Item(70, 11).Value = 7109
$ws.Cells.Item(70, 12).Value = 9666.667
$ws.Cells.Item(70, 13).Value = -6839
$ws.Cells.Item(70, 14).Value = -10206.667
$ws.Cells.Item(73, 8).Value = 8387.833
$ws.Cells.Item(73, 9).Value = 7109
$ws.Cells.Item(73, 10).Value = 9666.667
$ws.Cells.Item(73, 11).Value = 7109
$ws.Cells.Item(73, 12).Value = 9666.667
$ws.Cells.Item(73, 13).Value = -6173
$ws.Cells.Item(73, 14).Value = -11538.667
$ws.Cells.Item(80, 8).Value = 28950
$ws.Cells.Item(80, 9).Value = 8000
$ws.Cells.Item(80, 10).Value = 49900
$ws.Cells.Item(80, 11).Value = 8000
$ws.Cells.Item(80, 12).Value = 49900
$ws.Cells.Item(80, 13).Value = -7002
$ws.Cells.Item(80, 14).Value = -51896
$ws.Cells.Item(83, 8).Value = 28950
$ws.Cells.Item(83, 9).Value = 8000
$ws.Cells.Item(83, 10).Value = 49900
$ws.Cells.Item(83, 11).Value = 40000
$ws.Cells.Item(83, 12).Value = 249500
$ws.Cells.Item(83, 13).Value = -35008
$ws.Cells.Item(83, 14).Value = -259484
$ws.Cells.Item(132, 8).Value = 2727.3845
$ws.Cells.Item(132, 9).Value = 2306.2222
$ws.Cells.Item(132, 11).Value = 6918.6666
$ws.Cells.Item(132, 13).Value = -4388.6666
$ws.Cells.Item(135, 8).Value = 75000
$ws.Cells.Item(135, 10).Value = 75000
$ws.Cells.Item(135, 12).Value = 75000
$ws.Cells.Item(135, 14).Value = -85140

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(56, 8).Value = 21889
$ws.Cells.Item(56, 9).Value = 12899.667
$ws.Cells.Item(56, 11).Value = 12899.667
$ws.Cells.Item(56, 13).Value = -12208.667
$ws.Cells.Item(76, 8).Value = 3000
$ws.Cells.Item(76, 10).Value = 3000
$ws.Cells.Item(76, 12).Value = 3000
$ws.Cells.Item(76, 14).Value = -3676
$ws.Cells.Item(79, 8).Value = 3000
$ws.Cells.Item(79, 10).Value = 3000
$ws.Cells.Item(79, 12).Value = 3000
$ws.Cells.Item(79, 14).Value = -5340
$ws.Cells.Item(82, 8).Value = 1380.2
$ws.Cells.Item(82, 10).Value = 1950
$ws.Cells.Item(82, 12).Value = 1950
$ws.Cells.Item(82, 14).Value = -2672
$ws.Cells.Item(85, 8).Value = 1380.2
$ws.Cells.Item(85, 10).Value = 1950
$ws.Cells.Item(85, 12).Value = 1950
$ws.Cells.Item(85, 14).Value = -4446
$ws.Cells.Item(93, 8).Value = 1050
$ws.Cells.Item(93, 9).Value = 1050
$ws.Cells.Item(93, 11).Value = 1050
$ws.Cells.Item(93, 13).Value = 198

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2306.25
$ws.Cells.Item(96, 9).Value = 1996.4
$ws.Cells.Item(96, 10).Value = 2822.6667
$ws.Cells.Item(96, 11).Value = 1996.4
$ws.Cells.Item(96, 12).Value = 2822.6667
$ws.Cells.Item(96, 13).Value = -623.4000000000001
